# Update the cached date/time placeholder text from 09.03.22 to 10.03.22
# on the NotesMaster, the SlideMaster and every slide layout (these are
# auto-updating "last printed/edited" date fields that PowerPoint re-caches
# whenever the deck is saved on a different day), and fix a typo
# ("disctings" -> "distinct") in the body text of slide 8.

$p = $ppt.ActivePresentation

$oldDate = "09.03.22"
$newDate = "10.03.22"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if ($isDatePlaceholder -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# --- NotesMaster date placeholder ---
Update-DatePlaceholder $p.NotesMaster.Shapes

# --- SlideMaster date placeholder ---
$slideMaster = $p.SlideMaster
Update-DatePlaceholder $slideMaster.Shapes

# --- Every slide layout's date placeholder ---
$layouts = $slideMaster.CustomLayouts
for ($l = 1; $l -le $layouts.Count; $l++) {
    Update-DatePlaceholder $layouts.Item($l).Shapes
}

# --- Slide 8, "Textfeld 4": fix "disctings" -> "distinct" ---
$slide = $p.Slides.Item(8)
$shape = $slide.Shapes.Item("Textfeld 4")
$textRange = $shape.TextFrame.TextRange
$paragraph = $textRange.Paragraphs(11, 1)

$oldRun = "Simple comparision: gridsearch on 6 parameters with e.g. 10 disctings steps for each parameter means to try on "
$newRun = "Simple comparision: gridsearch on 6 parameters with e.g. 10 distinct steps for each parameter means to try on "

$runRange = $textRange.Characters($paragraph.Start, $oldRun.Length)
if ($runRange.Text -eq $oldRun) {
    $runRange.Text = $newRun
}
